# Add a new worksheet "Sheet2" after the existing "Employees" sheet.
# It reproduces the "Employees" table shifted one column right / two rows
# down (B3:E4), plus a brand-new third data row (B5:E5) that introduces two
# new shared strings ("Die" / "Director").
#
# We build Sheet2 by duplicating "Employees" (Worksheets.Copy) rather than
# Worksheets.Add(), because Copy faithfully carries over the existing
# per-cell formatting (cell styles) *and* the worksheet's drawing
# relationship in one step - matching the source table's look exactly.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Duplicate "Employees" and place the copy as the last sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws1.Copy([System.Reflection.Missing]::Value, $lastSheet)
$ws2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2.Name = "Sheet2"

# Shift the copied A1:D2 block to B3:E4, then wipe the now-empty A1:D2
# region (both values and the formatting that Cut leaves behind).
$ws2.Range("A1:D2").Cut($ws2.Range("B3"))
$ws2.Range("A1:D2").Clear()

# Append a third row (B5:E5), cloning row 4's formatting, then overwrite
# with the new data - including the two new strings "Die" and "Director".
$ws2.Range("B4:E4").Copy($ws2.Range("B5"))
$ws2.Range("B5").Value = 3.0
$ws2.Range("C5").Value = "Pew"
$ws2.Range("D5").Value = "Die"
$ws2.Range("E5").Value = "Director"

# Restore the original active sheet/selection.
$ws1.Select()
